$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 154.94118
$ws.Range("I33").Value = 172.07692
$ws.Range("K33").Value = 172.07692
$ws.Range("M33").Value = 56.92308
$ws.Range("H39").Value = 16.6
$ws.Range("J39").Value = 29.5
$ws.Range("L39").Value = 88.5
$ws.Range("N39").Value = -680.5
$ws.Range("H62").Value = 4666
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4999
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4999
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6247
$ws.Range("H65").Value = 4666
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4999
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 24995
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -31235
$ws.Range("H98").Value = 1683.1333
$ws.Range("I98").Value = 1341.72
$ws.Range("K98").Value = 1341.72
$ws.Range("M98").Value = 156.28
$ws.Range("H122").Value = 1683.1333
$ws.Range("I122").Value = 1341.72
$ws.Range("K122").Value = 4025.16
$ws.Range("M122").Value = -1575.16
$ws.Range("H132").Value = 5832.87
$ws.Range("I132").Value = 1414.7018
$ws.Range("J132").Value = 11689.512
$ws.Range("K132").Value = 4244.1054
$ws.Range("L132").Value = 35068.536
$ws.Range("M132").Value = -1714.1054
$ws.Range("N132").Value = -40128.536
$ws.Range("H135").Value = 8488.6
$ws.Range("I135").Value = 2462
$ws.Range("K135").Value = 22158
$ws.Range("M135").Value = -19623
$ws.Range("H137").Value = 5879.25
$ws.Range("I137").Value = 2500
$ws.Range("J137").Value = 7005.6665
$ws.Range("K137").Value = 7500
$ws.Range("L137").Value = 21016.9995
$ws.Range("M137").Value = -4950
$ws.Range("N137").Value = -26116.9995
$ws.Range("H140").Value = 71816.5
$ws.Range("J140").Value = 69541.25
$ws.Range("L140").Value = 69541.25
$ws.Range("N140").Value = -79901.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4492287
$ws.Range("I2").Value = 5823050
$ws.Range("J2").Value = 499999
$ws.Range("K2").Value = 5823050
$ws.Range("L2").Value = 499999
$ws.Range("M2").Value = -5822937
$ws.Range("N2").Value = -500225
$ws.Range("H32").Value = 4107.1333
$ws.Range("I32").Value = 2334.0652
$ws.Range("J32").Value = 9932.929
$ws.Range("K32").Value = 2334.0652
$ws.Range("L32").Value = 9932.929
$ws.Range("M32").Value = -2047.0652
$ws.Range("N32").Value = -10506.929
$ws.Range("H61").Value = 9473.209000000001
$ws.Range("I61").Value = 7366.909
$ws.Range("K61").Value = 7366.909
$ws.Range("M61").Value = -7154.909
$ws.Range("H74").Value = 1571.4286
$ws.Range("H77").Value = 1571.4286
$ws.Range("H116").Value = 4492287
$ws.Range("I116").Value = 5823050
$ws.Range("J116").Value = 499999
$ws.Range("K116").Value = 5823050
$ws.Range("L116").Value = 499999
$ws.Range("M116").Value = -5820756
$ws.Range("N116").Value = -504587
$ws.Range("H132").Value = 11531.65
$ws.Range("I132").Value = 13296.174
$ws.Range("K132").Value = 39888.522
$ws.Range("M132").Value = -37358.522
$ws.Range("H136").Value = 9473.209000000001
$ws.Range("I136").Value = 7366.909
$ws.Range("K136").Value = 22100.727
$ws.Range("M136").Value = -19550.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4492287
$ws.Range("I3").Value = 5823050
$ws.Range("J3").Value = 499999
$ws.Range("K3").Value = 5823050
$ws.Range("L3").Value = 499999
$ws.Range("M3").Value = -5822936
$ws.Range("N3").Value = -500227
$ws.Range("H86").Value = 2832.9656
$ws.Range("I86").Value = 2588.1875
$ws.Range("J86").Value = 3134.2307
$ws.Range("K86").Value = 2588.1875
$ws.Range("L86").Value = 3134.2307
$ws.Range("M86").Value = -1465.1875
$ws.Range("N86").Value = -5380.2307
$ws.Range("H89").Value = 2832.9656
$ws.Range("I89").Value = 2588.1875
$ws.Range("J89").Value = 3134.2307
$ws.Range("K89").Value = 12940.9375
$ws.Range("L89").Value = 15671.1535
$ws.Range("M89").Value = -7324.9375
$ws.Range("N89").Value = -26903.1535
$ws.Range("H134").Value = 1425.6154
$ws.Range("I134").Value = 1168.8334
$ws.Range("K134").Value = 3506.5002
$ws.Range("M134").Value = -971.5001999999999
$ws.Range("H140").Value = 144647.2
$ws.Range("J140").Value = 144647.2
$ws.Range("L140").Value = 144647.2
$ws.Range("N140").Value = -155007.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1306.8
$ws.Range("J16").Value = 1765
$ws.Range("L16").Value = 1765
$ws.Range("N16").Value = -2339
$ws.Range("H31").Value = 3767.3274
$ws.Range("I31").Value = 3063.1086
$ws.Range("K31").Value = 3063.1086
$ws.Range("M31").Value = -2768.1086
$ws.Range("H34").Value = 3767.3274
$ws.Range("I34").Value = 3063.1086
$ws.Range("K34").Value = 3063.1086
$ws.Range("M34").Value = -2861.1086
$ws.Range("H58").Value = 385921.56
$ws.Range("I58").Value = 527449.9399999999
$ws.Range("K58").Value = 527449.9399999999
$ws.Range("M58").Value = -527246.9399999999
$ws.Range("H99").Value = 17721.066
$ws.Range("I99").Value = 21865.375
$ws.Range("J99").Value = 12984.714
$ws.Range("K99").Value = 21865.375
$ws.Range("L99").Value = 12984.714
$ws.Range("M99").Value = -20367.375
$ws.Range("N99").Value = -15980.714
$ws.Range("H107").Value = 588202.1
$ws.Range("I107").Value = 1069788.2
$ws.Range("K107").Value = 1069788.2
$ws.Range("M107").Value = -1067868.2
$ws.Range("H113").Value = 1306.8
$ws.Range("J113").Value = 1765
$ws.Range("L113").Value = 1765
$ws.Range("N113").Value = -6105
$ws.Range("H126").Value = 17721.066
$ws.Range("I126").Value = 21865.375
$ws.Range("J126").Value = 12984.714
$ws.Range("K126").Value = 65596.125
$ws.Range("L126").Value = 38954.142
$ws.Range("M126").Value = -63126.125
$ws.Range("N126").Value = -43894.142
$ws.Range("H132").Value = 3701.9
$ws.Range("I132").Value = 3377.5
$ws.Range("K132").Value = 10132.5
$ws.Range("M132").Value = -7602.5
$ws.Range("H134").Value = 1913.9166
$ws.Range("I134").Value = 1666.8334
$ws.Range("K134").Value = 5000.5002
$ws.Range("M134").Value = -2465.5002
$ws.Range("H136").Value = 385921.56
$ws.Range("I136").Value = 527449.9399999999
$ws.Range("K136").Value = 1582349.82
$ws.Range("M136").Value = -1579799.82
$ws.Range("H141").Value = 85590.64999999999
$ws.Range("I141").Value = 41999.5
$ws.Range("J141").Value = 91402.8
$ws.Range("K141").Value = 41999.5
$ws.Range("L141").Value = 91402.8
$ws.Range("M141").Value = -36819.5
$ws.Range("N141").Value = -101762.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9003.666999999999
$ws.Range("I3").Value = 3671.5557
$ws.Range("K3").Value = 11014.6671
$ws.Range("M3").Value = -10902.6671
$ws.Range("H5").Value = 1085.5
$ws.Range("I5").Value = 926.2857
$ws.Range("K5").Value = 2778.8571
$ws.Range("M5").Value = -2666.8571
$ws.Range("H97").Value = 402
$ws.Range("I97").Value = 402
$ws.Range("K97").Value = 1206
$ws.Range("M97").Value = -710
$ws.Range("H135").Value = 1085.5
$ws.Range("I135").Value = 926.2857
$ws.Range("K135").Value = 8336.5713
$ws.Range("M135").Value = -5801.5713
$ws.Range("H137").Value = 8475385
$ws.Range("I137").Value = 3335.1667
$ws.Range("J137").Value = 15737143
$ws.Range("K137").Value = 10005.5001
$ws.Range("L137").Value = 47211429
$ws.Range("M137").Value = -4905.500100000001
$ws.Range("N137").Value = -47221629

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1164.4736
$ws.Range("I93").Value = 1139.8462
$ws.Range("K93").Value = 1139.8462
$ws.Range("M93").Value = 108.1538
$ws.Range("H111").Value = 75193.5
$ws.Range("J111").Value = 75193.5
$ws.Range("L111").Value = 75193.5
$ws.Range("N111").Value = -83373.5
$ws.Range("H136").Value = 6500.25
$ws.Range("I136").Value = 5833.8335
$ws.Range("K136").Value = 17501.5005
$ws.Range("M136").Value = -14951.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 60000
$ws.Range("J49").Value = 60000
$ws.Range("L49").Value = 60000
$ws.Range("N49").Value = -60460
$ws.Range("H81").Value = 4634640
$ws.Range("J81").Value = 10419666
$ws.Range("L81").Value = 20839332
$ws.Range("N81").Value = -20841454
$ws.Range("H84").Value = 4634640
$ws.Range("J84").Value = 10419666
$ws.Range("L84").Value = 104196660
$ws.Range("N84").Value = -104207268
$ws.Range("H132").Value = 22229730
$ws.Range("I132").Value = 3474071.5
$ws.Range("K132").Value = 10422214.5
$ws.Range("M132").Value = -10419684.5
$ws.Range("H136").Value = 8038.2065
$ws.Range("I136").Value = 3264.375
$ws.Range("J136").Value = 9043.224
$ws.Range("K136").Value = 9793.125
$ws.Range("L136").Value = 27129.672
$ws.Range("M136").Value = -7243.125
$ws.Range("N136").Value = -32229.672
